# Applies hybrid bold + color (2C3E50) highlighting to quantitative impact
# metrics (percentages / dollar amounts) inside specific bullet paragraphs,
# matching the target XML diff. Each paragraph is located first by its full
# original (unique) text via Find.Execute, which also narrows the returned
# Range to the paragraph bounds; the bold "metric" substrings are then
# located and formatted within a Range clamped to those same bounds so nothing
# outside the paragraph is ever touched.

$d = $word.ActiveDocument

# Hex color 2C3E50 -> Word/VBA BGR-packed long (R + G*256 + B*65536)
$metricColor = 5258796

function Highlight-Metrics($paraText, $metrics) {
    $rng = $d.Content
    $found = $rng.Find.Execute($paraText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $paraText"
        return
    }
    $pStart = $rng.Start
    $pEnd = $rng.End
    foreach ($m in $metrics) {
        $sub = $d.Range($pStart, $pEnd)
        $f = $sub.Find.Execute($m, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($f) {
            $sub.Font.Bold = 1
            $sub.Font.Color = $metricColor
        } else {
            Write-Output "METRIC NOT FOUND: $m in $paraText"
        }
    }
}

# 1. Siege Analytics - race coding errors bullet
Highlight-Metrics "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" @("23%", "64%")

# 2. Siege Analytics - advanced sampling methods bullet
Highlight-Metrics "• Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes" @("±4.2%", "±2.1%", "71%", "87%")

# 3. Siege Analytics - trigonometric algorithm bullet
Highlight-Metrics "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis" @("73.5%", "`$4.7M")

# 4. Siege Analytics - FEC analysis systems bullet
Highlight-Metrics "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion" @("`$2")

# 5. Key Achievements - algorithmic innovation bullet
Highlight-Metrics "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%" @("73.5%")

# 6. Key Achievements - $4.7M savings bullet
Highlight-Metrics "• `$4.7M savings enabled nonprofit access" @("`$4.7M")

# 7. Key Achievements - 178% accuracy improvement bullet
Highlight-Metrics "• 178% accuracy improvement in racial classification algorithms" @("178%")
